$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unused "KS_iMax" column (column H). Deleting the entire
# column shifts the later columns (insilico_call, inSilico_AD) one
# position to the left, matching the cleaned-up layout, and also drops
# the now-unreferenced "KS_iMax" shared string.
$ws.Range("H1").EntireColumn.Delete()

# Reset the selection back to the top-left cell (UI cleanup).
$ws.Range("A1").Select()
